# Apply the "2025-09-19 06:26 JST" scrape refresh to the ランサーズ sheet:
#  - new fetch timestamp for every kept row
#  - rows 3-7 replaced with new listing data
#  - rows 8-20 removed (only 6 listings remain this run)
#  - column B/D widths tweaked
#  - hyperlinks rebuilt so only F2:F7 carry links

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-19 06:26:53"

# --- new row data for rows 3-7 (row 2 keeps its old content, only the
#     timestamp changes) ---
$rows = @(
    @{ B = "システムの開発補助や運営サポート【フルリモート×長期】";
       C = "システム開発";
       D = "200,000 円 ~ 300,000 円 / 固定";
       E = "期限情報なし";
       F = "https://www.lancers.jp/work/detail/5323359";
       G = 83;
       H = "◆開発" },
    @{ B = "初回 【フィンテック/ブリッジ】金融資産管理システムの要件定義/仕様伝達/進捗管理(日/英|フルリモート)";
       C = "システム開発";
       D = "300,000 円 ~ 500,000 円 / 固定";
       E = "期限情報なし";
       F = "https://www.lancers.jp/work/detail/5396502";
       G = 60;
       H = "◇管理" },
    @{ B = "【フィンテック/QA】海外の金融資産管理システムのテスト設計・品質保証";
       C = "システム開発";
       D = "200,000 円 ~ 300,000 円 / 固定";
       E = "期限情報なし";
       F = "https://www.lancers.jp/work/detail/5396510";
       G = 53;
       H = "◇管理" },
    @{ B = "【急募】PHP・Lalavelでの既存プログラム改修依頼";
       C = "システム開発";
       D = "100,000 円 ~ 200,000 円 / 固定";
       E = "期限情報なし";
       F = "https://www.lancers.jp/work/detail/5396563";
       G = 33;
       H = "○PHP" },
    @{ B = "【急募】WordPress記事をCoopelで自動投稿設定できる方を探しています!";
       C = "システム開発";
       D = "5,000 円 ~ 10,000 円 / 固定";
       E = "期限情報なし";
       F = "https://www.lancers.jp/work/detail/5396609";
       G = 25;
       H = "○WordPress" }
)

# Remove old hyperlinks first (they get rebuilt below to match the
# surviving rows only).
$ws.Hyperlinks.Delete()

# Refresh the "取得日時" column for the rows that remain (2-7).
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Write the new listing content into rows 3-7.
$r = 3
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# Drop the now-stale rows 8-20 (was 20 listings, now only 6).
$ws.Range("A8:H20").EntireRow.Delete()

# Rebuild the hyperlinks for the surviving URL column (F2:F7).
for ($r = 2; $r -le 7; $r++) {
    $url = $ws.Cells.Item($r, 6).Value()
    $ws.Hyperlinks.Add($ws.Range("F" + $r), $url)
}

# Column width adjustments (B: 46 -> 55, D: 32 -> 28).
# Note: the host round-trips ColumnWidth through a pixel conversion that
# adds ~0.83 back on save, so dial in 0.83 less than the target stored
# width to land exactly on 55 / 28 in the saved OOXML.
$ws.Columns.Item(2).ColumnWidth = 54.17
$ws.Columns.Item(4).ColumnWidth = 27.17
